$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1740.25
$ws.Range("I40").Value = 1663
$ws.Range("J40").Value = 1817.5
$ws.Range("K40").Value = 1663
$ws.Range("L40").Value = 1817.5
$ws.Range("M40").Value = -1488
$ws.Range("N40").Value = -2167.5

$ws.Range("H62").Value = 3880
$ws.Range("I62").Value = 4033.3333
$ws.Range("J62").Value = 3650
$ws.Range("K62").Value = 4033.3333
$ws.Range("L62").Value = 3650
$ws.Range("M62").Value = -3409.3333
$ws.Range("N62").Value = -4898

$ws.Range("H65").Value = 3880
$ws.Range("I65").Value = 4033.3333
$ws.Range("J65").Value = 3650
$ws.Range("K65").Value = 20166.6665
$ws.Range("L65").Value = 18250
$ws.Range("M65").Value = -17046.6665
$ws.Range("N65").Value = -24490

$ws.Range("H82").Value = 1056.8
$ws.Range("I82").Value = 1056.8
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3170.4
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2764.4
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 1056.8
$ws.Range("I85").Value = 1056.8
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3170.4
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1766.4
$ws.Range("N85").ClearContents()

$ws.Range("H111").Value = 10225
$ws.Range("I111").Value = 10300
$ws.Range("J111").Value = 10000
$ws.Range("K111").Value = 30900
$ws.Range("L111").Value = 30000
$ws.Range("M111").Value = -27833
$ws.Range("N111").Value = -36134

$ws.Range("H123").Value = 43899.668
$ws.Range("J123").Value = 43899.668
$ws.Range("L123").Value = 43899.668
$ws.Range("N123").Value = -53699.668

$ws.Range("H125").Value = 761.61536
$ws.Range("I125").Value = 560
$ws.Range("J125").Value = 887.625
$ws.Range("K125").Value = 5040
$ws.Range("L125").Value = 7988.625
$ws.Range("M125").Value = -2580
$ws.Range("N125").Value = -12908.625

$ws.Range("H132").Value = 1436.2046
$ws.Range("I132").Value = 1392.325
$ws.Range("J132").Value = 1875
$ws.Range("K132").Value = 4176.975
$ws.Range("L132").Value = 5625
$ws.Range("M132").Value = -1646.975
$ws.Range("N132").Value = -10685

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2044.5217
$ws.Range("I45").Value = 1923.5555
$ws.Range("J45").Value = 2480
$ws.Range("K45").Value = 1923.5555
$ws.Range("L45").Value = 2480
$ws.Range("M45").Value = -1546.5555
$ws.Range("N45").Value = -3234

$ws.Range("H54").Value = 18000
$ws.Range("J54").Value = 18000
$ws.Range("L54").Value = 18000
$ws.Range("N54").Value = -19538

$ws.Range("H61").Value = 4233.237
$ws.Range("I61").Value = 3494.6924
$ws.Range("J61").Value = 5833.4165
$ws.Range("K61").Value = 3494.6924
$ws.Range("L61").Value = 5833.4165
$ws.Range("M61").Value = -3282.6924
$ws.Range("N61").Value = -6257.4165

$ws.Range("H63").Value = 4187
$ws.Range("I63").Value = 4187
$ws.Range("K63").Value = 4187
$ws.Range("M63").Value = -3501

$ws.Range("H66").Value = 4187
$ws.Range("I66").Value = 4187
$ws.Range("K66").Value = 20935
$ws.Range("M66").Value = -17503

$ws.Range("H74").Value = 1620.8235
$ws.Range("I74").Value = 1690
$ws.Range("J74").Value = 514
$ws.Range("K74").Value = 1690
$ws.Range("L74").Value = 514
$ws.Range("M74").Value = -816
$ws.Range("N74").Value = -2262

$ws.Range("H77").Value = 1620.8235
$ws.Range("I77").Value = 1690
$ws.Range("J77").Value = 514
$ws.Range("K77").Value = 8450
$ws.Range("L77").Value = 2570
$ws.Range("M77").Value = -4082
$ws.Range("N77").Value = -11306

$ws.Range("H123").Value = 52700
$ws.Range("J123").Value = 52700
$ws.Range("L123").Value = 52700
$ws.Range("N123").Value = -62500

$ws.Range("H136").Value = 4233.237
$ws.Range("I136").Value = 3494.6924
$ws.Range("J136").Value = 5833.4165
$ws.Range("K136").Value = 10484.0772
$ws.Range("L136").Value = 17500.2495
$ws.Range("M136").Value = -7934.0772
$ws.Range("N136").Value = -22600.2495

$ws.Range("H140").Value = 51084.8
$ws.Range("J140").Value = 51084.8
$ws.Range("L140").Value = 51084.8
$ws.Range("N140").Value = -61444.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 47305
$ws.Range("J140").Value = 47305
$ws.Range("L140").Value = 47305
$ws.Range("N140").Value = -57665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2261.1458
$ws.Range("I31").Value = 1821.0834
$ws.Range("J31").Value = 2701.2083
$ws.Range("K31").Value = 1821.0834
$ws.Range("L31").Value = 2701.2083
$ws.Range("M31").Value = -1526.0834
$ws.Range("N31").Value = -3291.2083

$ws.Range("H34").Value = 2261.1458
$ws.Range("I34").Value = 1821.0834
$ws.Range("J34").Value = 2701.2083
$ws.Range("K34").Value = 1821.0834
$ws.Range("L34").Value = 2701.2083
$ws.Range("M34").Value = -1619.0834
$ws.Range("N34").Value = -3105.2083

$ws.Range("H41").Value = 8500
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H43").Value = 38000
$ws.Range("J43").Value = 38000
$ws.Range("L43").Value = 38000
$ws.Range("N43").Value = -38368

$ws.Range("H47").Value = 16658.375
$ws.Range("I47").Value = 9000
$ws.Range("J47").Value = 17752.428
$ws.Range("K47").Value = 9000
$ws.Range("L47").Value = 17752.428
$ws.Range("M47").Value = -8434
$ws.Range("N47").Value = -18884.428

$ws.Range("H53").Value = 34074.75
$ws.Range("J53").Value = 34074.75
$ws.Range("L53").Value = 34074.75
$ws.Range("N53").Value = -35288.75

$ws.Range("H58").Value = 2600476.2
$ws.Range("I58").Value = 6995249
$ws.Range("K58").Value = 6995249
$ws.Range("M58").Value = -6995046

$ws.Range("H101").Value = 38000
$ws.Range("J101").Value = 38000
$ws.Range("L101").Value = 38000
$ws.Range("N101").Value = -44490

$ws.Range("H132").Value = 2018.4584
$ws.Range("I132").Value = 1544.2858
$ws.Range("J132").Value = 5337.6665
$ws.Range("K132").Value = 4632.857400000001
$ws.Range("L132").Value = 16012.9995
$ws.Range("M132").Value = -2102.857400000001
$ws.Range("N132").Value = -21072.9995

$ws.Range("H134").Value = 2633.9216
$ws.Range("I134").Value = 1863.9656
$ws.Range("J134").Value = 3648.8635
$ws.Range("K134").Value = 5591.8968
$ws.Range("L134").Value = 10946.5905
$ws.Range("M134").Value = -3056.8968
$ws.Range("N134").Value = -16016.5905

$ws.Range("H136").Value = 2600476.2
$ws.Range("I136").Value = 6995249
$ws.Range("K136").Value = 20985747
$ws.Range("M136").Value = -20983197

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 4016
$ws.Range("I82").Value = 1078.6
$ws.Range("J82").Value = 4995.1333
$ws.Range("K82").Value = 3235.8
$ws.Range("L82").Value = 14985.3999
$ws.Range("M82").Value = -2829.8
$ws.Range("N82").Value = -15797.3999

$ws.Range("H85").Value = 4016
$ws.Range("I85").Value = 1078.6
$ws.Range("J85").Value = 4995.1333
$ws.Range("K85").Value = 3235.8
$ws.Range("L85").Value = 14985.3999
$ws.Range("M85").Value = -1831.8
$ws.Range("N85").Value = -17793.3999

$ws.Range("H121").Value = 826.25
$ws.Range("I121").Value = 140
$ws.Range("J121").Value = 1708.5714
$ws.Range("K121").Value = 420
$ws.Range("L121").Value = 5125.7142
$ws.Range("M121").Value = 890
$ws.Range("N121").Value = -7745.7142

$ws.Range("H131").Value = 879.1539
$ws.Range("I131").Value = 655
$ws.Range("J131").Value = 953.87177
$ws.Range("K131").Value = 1965
$ws.Range("L131").Value = 2861.61531
$ws.Range("M131").Value = 3075
$ws.Range("N131").Value = -12941.61531

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 816.06665
$ws.Range("I107").Value = 597.9
$ws.Range("J107").Value = 1252.4
$ws.Range("K107").Value = 597.9
$ws.Range("L107").Value = 1252.4
$ws.Range("M107").Value = 1322.1
$ws.Range("N107").Value = -5092.4

$ws.Range("H113").Value = 2845.3845
$ws.Range("I113").Value = 2439.3
$ws.Range("J113").Value = 4199
$ws.Range("K113").Value = 2439.3
$ws.Range("L113").Value = 4199
$ws.Range("M113").Value = -269.3000000000002
$ws.Range("N113").Value = -8539

$ws.Range("H126").Value = 2685.7856
$ws.Range("I126").Value = 1851
$ws.Range("J126").Value = 3311.875
$ws.Range("K126").Value = 5553
$ws.Range("L126").Value = 9935.625
$ws.Range("M126").Value = -3083
$ws.Range("N126").Value = -14875.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5253.0303
$ws.Range("I122").Value = 5048.857
$ws.Range("J122").Value = 6396.4
$ws.Range("K122").Value = 15146.571
$ws.Range("L122").Value = 19189.2
$ws.Range("M122").Value = -12696.571
$ws.Range("N122").Value = -24089.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H136").Value = 2657.8386
$ws.Range("I136").Value = 2078.7188
$ws.Range("J136").Value = 3275.5667
$ws.Range("K136").Value = 6236.1564
$ws.Range("L136").Value = 9826.7001
$ws.Range("M136").Value = -3686.1564
$ws.Range("N136").Value = -14926.7001
